# edit.ps1 - apply cryptos.xlsx symbol/price update (commit: Mon Dec 12 19:25:39 UTC 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '274.76'
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.10'
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.201'
$ws.Range("D4").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.578'
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.516'
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '6.528'
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8230'
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1644'
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08231'
$ws.Range("D11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03424'
$ws.Range("D12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03140'
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09135'
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.771'
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001624'
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04697'
$ws.Range("D17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006453'
$ws.Range("D18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006137'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '18HotbitTokenHTB'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.722'
$ws.Range("D22").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3276'
$ws.Range("D25").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04743'
$ws.Range("D40").ClearFormats()

$ws.Range("B41").Value = 'CEJI'

$ws.Range("C41").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005503'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = '40CEJICEJI'

$ws.Range("B42").Value = 'KickToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007029'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '41KickTokenKICKBestin24h'

$ws.Range("B43").Value = 'BKEXToken'

$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1103'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '42BKEXTokenBKK'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01028'
$ws.Range("D44").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00001901'
$ws.Range("D49").ClearFormats()

